$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) A7: add the new "abstract" paragraph (plain run, keeps the existing
#    bordered / wrap / vertical-centered style already used by column A).
# ---------------------------------------------------------------------------
$paragraph = 'Las células solares orgánicas nos muestran un camino interesante hacia el uso de energías renovables ecológicas y amigables. Eso nos ayudará a mitigar el efecto o huella de carbono. Se están cuestionando formas eficientes de convertir la energía solar en electricidad, como el uso de materiales, buscando las mejores propiedades que permitan una óptima conversión energética. Este trabajo explora el uso de técnicas de aprendizaje automático (ML) para ayudar a optimizar propiedades moleculares como el orbital molecular de alta ocupación (HOMO) y las energías de orbital molecular desocupado más bajo (LUMO), así como el cálculo y calibración de la eficiencia de conversión de potencia (PCE). ) con el ánimo de buscar grandes candidatas a moléculas orgánicas para su uso como sistemas donante-receptor en células solares. En particular, probamos una calibración del proceso gaussiano como un modelo ML en un conjunto de moléculas reportadas en la literatura [1] y discutimos algunos aspectos tanto de las propiedades químicas como de la ventaja de usar ML'

$a7 = $ws.Range("A7")
$a7.Value = $paragraph

# ---------------------------------------------------------------------------
# 2) A8: add the "Keywords:" line. "Keywords:" stays bold+italic Times New
#    Roman (the cell's own font/style), the rest of the sentence is plain
#    (non-bold/non-italic) Times New Roman applied to that run only.
# ---------------------------------------------------------------------------
$kwLabel = 'Keywords:'
$kwRest  = ' Organic Solar Cell, Small Molecules, Machine Learning, Computational Chemistry, Quantum Systems.'
$kwFull  = $kwLabel + $kwRest

$a8 = $ws.Range("A8")
$a8.Value = $kwFull

# Cell-level (and therefore "Keywords:" run) formatting: bold italic Times
# New Roman, no border, justified horizontally, vertically centered, no wrap.
$a8.Font.Name = "Times New Roman"
$a8.Font.Family = 1
$a8.Font.Bold = $true
$a8.Font.Italic = $true
$a8.Borders.LineStyle = -4142
$a8.WrapText = $false
$a8.HorizontalAlignment = -4130
$a8.VerticalAlignment = -4108

# Now re-format just the trailing sentence as a plain (non-bold, non-italic)
# Times New Roman run.
$restLen = $kwFull.Length - $kwLabel.Length
$restRange = $a8.Characters($kwLabel.Length + 1, $restLen)
$restRange.Font.Name = "Times New Roman"
$restRange.Font.Family = 1
$restRange.Font.Bold = $false
$restRange.Font.Italic = $false

# ---------------------------------------------------------------------------
# 3) Row heights (re-flowed once the new text above was added).
# ---------------------------------------------------------------------------
$ws.Rows(4).RowHeight = 75.75
$ws.Rows(7).RowHeight = 68.25
$ws.Rows(8).RowHeight = 30.75

# ---------------------------------------------------------------------------
# 4) Selection / active cell ends up on A8.
# ---------------------------------------------------------------------------
$a8.Select()
